$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.865.68"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.209.41"
$ws.Range("E3").Value = "  -1.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.18%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.39"
$ws.Range("E5").Value = "  -0.23%  "

# Row 6
$ws.Range("E6").Value = "  -1.57%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.64"
$ws.Range("E7").Value = "  -2.46%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.401"
$ws.Range("E9").Value = "  -1.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.14"
$ws.Range("E10").Value = "  -2.77%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0896"
$ws.Range("E11").Value = "  +1.23%  "

# Row 12
$ws.Range("E12").Value = "  -0.90%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.539.65"
$ws.Range("E13").Value = "  -1.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.40"
$ws.Range("E14").Value = "  -2.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.10"
$ws.Range("E15").Value = "  +1.43%  "

# Row 16
$ws.Range("E16").Value = "  -1.76%  "

# Row 17
$ws.Range("E17").Value = "  -1.69%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.230.20"
$ws.Range("E18").Value = "  -0.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.809.85"
$ws.Range("E19").Value = "  +0.36%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0922"
$ws.Range("E20").Value = "  +2.05%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.92"
$ws.Range("E21").Value = "  -2.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.06"
$ws.Range("E22").Value = "  -2.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.57"
$ws.Range("E23").Value = "  -2.15%  "

# Row 24
$ws.Range("E24").Value = "  -0.18%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  -0.36%  "

# Row 26
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.35"
$ws.Range("E26").Value = "  -2.79%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("E27").Value = "  -2.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.44"
$ws.Range("E28").Value = "  -0.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.139"
$ws.Range("E29").Value = "  -2.82%  "

# Row 30
$ws.Range("E30").Value = "  -3.14%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.74"
$ws.Range("E31").Value = "  -2.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.64"
$ws.Range("E32").Value = "  -5.36%  "

# Row 33
$ws.Range("E33").Value = "  -2.61%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  -3.01%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.59"
$ws.Range("E35").Value = "  -2.31%  "

# Row 36
$ws.Range("E36").Value = "  +2.50%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.33"
$ws.Range("E37").Value = "  -5.41%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.54"
$ws.Range("E38").Value = "  -8.02%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.32"
$ws.Range("E39").Value = "  -4.10%  "

# Row 40
$ws.Range("B40").Value = "BinanceUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.11%  "

# Row 41
$ws.Range("B41").Value = "TerraClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000237"
$ws.Range("E41").Value = "  -1.47%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0241"
$ws.Range("E42").Value = "  +0.96%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.54"
$ws.Range("E43").Value = "  -2.67%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0951"
$ws.Range("E44").Value = "  -2.10%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.19"
$ws.Range("E45").Value = "  -0.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "96.51"
$ws.Range("E46").Value = "  -4.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.451.00"
$ws.Range("E47").Value = "  -2.89%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.29"
$ws.Range("E48").Value = "  -12.56%  "

# Row 49
$ws.Range("E49").Value = "  -1.58%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.09"
$ws.Range("E50").Value = "  -3.40%  "

# Row 51
$ws.Range("E51").Value = "  -3.29%  "
